# lab_algo.xlsx — "Fixed plots, exponential search and matrix generations"
#
# The raw measurement table (BIN / EXP FIRST / STANDARD / EXP SECOND, rows
# 22:34 on Sheet1) was re-measured; the EXP2/EXP1 ratio column (J) is a
# shared formula (=H/D) so it recalculates on its own once D and H change.
# The scatter charts on the sheet just plot these same cells, so they track
# along with the edit automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New measurements for columns B (BIN), D (EXP FIRST), F (STANDARD),
# H (EXP SECOND), keyed by row. J (EXP2/EXP1 = H/D) is left to recalc.
$newData = @{
    22 = @{ B = 119;     D = 160;     F = 1524;  H = 132 }
    23 = @{ B = 410;     D = 370;     F = 2315;  H = 168 }
    24 = @{ B = 750;     D = 679;     F = 2703;  H = 211 }
    25 = @{ B = 1544;    D = 1359;    F = 2966;  H = 281 }
    26 = @{ B = 2966;    D = 2636;    F = 3209;  H = 414 }
    27 = @{ B = 6296;    D = 5549;    F = 3083;  H = 672 }
    28 = @{ B = 20322;   D = 13343;   F = 3346;  H = 1365 }
    29 = @{ B = 114817;  D = 31903;   F = 3780;  H = 3632 }
    30 = @{ B = 352328;  D = 74209;   F = 5433;  H = 10982 }
    31 = @{ B = 719388;  D = 142564;  F = 8279;  H = 35428 }
    32 = @{ B = 1405629; D = 299864;  F = 19223; H = 75942 }
    33 = @{ B = 3237251; D = 1335687; F = 55247; H = 211272 }
    34 = @{ B = 9543580; D = 3762314; F = 86582; H = 359028 }
}

foreach ($row in ($newData.Keys | Sort-Object)) {
    $vals = $newData[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B   # column B
    $ws.Cells.Item($row, 4).Value = $vals.D   # column D
    $ws.Cells.Item($row, 6).Value = $vals.F   # column F
    $ws.Cells.Item($row, 8).Value = $vals.H   # column H
}

# Sheet view: zoom + active selection moved.
$excel.ActiveWindow.Zoom = 114
$ws.Range("L21").Select()
